$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "LANGUAGE" -> "LANGUAGE_BUTTON" (shared string used by A3)
$ws.Range("A3").Value = "LANGUAGE_BUTTON"

# Re-apply the "Normal" cell style across the used range, which makes
# LibreOffice/Excel persist explicit font/alignment/protection "apply" flags
# on the cellXf instead of leaving them implicit.
$ws.Range("A1:D4").Style = "Normal"

# Move the active cell/selection to C3
$ws.Range("C3").Select()
